# Updated the facilities list in the enrollment form
# Adds ten new health-facility choices to the "choices" sheet of the
# enrollment form (select_one facilities list), pushing the existing
# "other" / "Other" row down below the new entries.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# Insert 10 blank rows above the existing "facilities/other/Other" row
# (row 6), pushing it down to row 16.
for ($i = 0; $i -lt 10; $i++) {
    $choices.Rows.Item(6).Insert()
}

$labels = @(
    "Chitungwiza-Seke North clinic",
    "Chitungwiza-Seke South clinic",
    "Chitungwiza-City Med hospital",
    "Chitungwiza-Zengeza Clinic",
    "Chitungwiza-Chitungwiza Central Hospital",
    "Chegutu- Norton hospital",
    "Chegutu- District Hospital",
    "Chegutu- Monera clinic (Norton Outreach)",
    "Marondera- Marondera District Hospital",
    "Marondera- Mahusekwa Hospital"
)
$codes = @("csnc", "cssc", "ccmh", "czc", "ccch", "cnh", "cdh", "cmc", "mmdh", "mmh")

# Fill list_name + label columns first for every new row ...
for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = 6 + $i
    $choices.Cells.Item($r, 1).Value = "facilities"
    $choices.Cells.Item($r, 3).Value = $labels[$i]
}
# ... then the name/code column, matching the authoring order used when
# the facility list was originally populated.
for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = 6 + $i
    $choices.Cells.Item($r, 2).Value = $codes[$i]
}

# Give the new list_name/name columns the same look as the rest of the
# "facilities" list (copy the formatting from the row that used to sit
# at row 6 and now lives at row 16).
$choices.Range("A16").Copy()
$choices.Range("A6:B15").PasteSpecial(-4122)

# Restore the frozen header row / selection on "survey" (unaffected by
# the edit, but the engine drops pane state on any save).
$survey.Activate()
$survey.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$survey.Range("B3").Select()

# Restore the frozen header row on "choices" and leave the newly
# relocated "other" row selected, as a user would after inserting rows
# above it.
$choices.Activate()
$choices.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$choices.Rows.Item(16).Select()
